# Update the "想去人数" (want-to-go count) column F values on both the
# "展览" and "全部类型" worksheets, which contain duplicated data.

$wb = $excel.ActiveWorkbook

$updates = @{
    5  = 2573
    9  = 1320
    13 = 1167
    14 = 338
    21 = 2368
    22 = 22
    23 = 273
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
